$d = $word.ActiveDocument

$replacements = @(
    @{old = "60×71=4260"; new = "80×68=5440"},
    @{old = "49×61=2989"; new = "77×88=6776"},
    @{old = "81×46=3726"; new = "50×79=3950"},
    @{old = "79×48=3792"; new = "75×58=4350"},
    @{old = "37×66=2442"; new = "87×12=1044"},
    @{old = "25×65=1625"; new = "62×80=4960"},
    @{old = "48×41=1968"; new = "91×70=6370"},
    @{old = "88×26=2288"; new = "11×92=1012"},
    @{old = "70×37=2590"; new = "75×75=5625"},
    @{old = "96×90=8640"; new = "21×61=1281"},
    @{old = "28×14=392";  new = "83×62=5146"},
    @{old = "86×39=3354"; new = "65×42=2730"},
    @{old = "98×43=4214"; new = "84×32=2688"},
    @{old = "16×23=368";  new = "19×80=1520"},
    @{old = "17×52=884";  new = "96×60=5760"},
    @{old = "82×29=2378"; new = "37×14=518"},
    @{old = "57×45=2565"; new = "23×77=1771"},
    @{old = "61×61=3721"; new = "56×37=2072"},
    @{old = "83×90=7470"; new = "24×20=480"},
    @{old = "11×13=143";  new = "14×37=518"},
    @{old = "22×90=1980"; new = "98×51=4998"},
    @{old = "12×63=756";  new = "19×29=551"},
    @{old = "20×50=1000"; new = "89×19=1691"},
    @{old = "52×50=2600"; new = "26×81=2106"},
    @{old = "77×17=1309"; new = "67×24=1608"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Output "Done applying replacements"
